$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns (I0, IF), matching the style of the
# existing header row (column H = "IP"). Copy formats only (xlPasteFormats)
# so the underlying value isn't clobbered by the paste.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I53 and J2:J53 (one entry per row 2..53)
$iVals = @(3,2,5,6,7,1,7,7,7,5,6,5,6,8,6,8,7,6,6,6,7,7,6,8,8,7,6,8,6,8,7,6,7,8,7,7,1,6,8,6,1,1,7,5,6,6,6,7,5,4,8,4)
$jVals = @(4,4,7,8,7,2,7,7,8,7,6,6,6,8,6,8,7,7,7,6,7,7,7,8,8,8,7,8,6,8,7,6,7,9,8,7,4,7,9,8,5,4,8,7,7,8,6,7,6,5,8,4)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
